# Replace each arithmetic-problem text in the table cells with its
# regenerated counterpart. Every "old" string in the document is unique,
# so a straightforward whole-word Find/Replace per pair is safe and
# order-independent (no "new" value collides with any "old" value).
$d = $word.ActiveDocument
$d.Content.Find.Execute("77-26=51", $true, $false, $false, $false, $false, $true, 1, $false, "42-3=39", 2) | Out-Null
$d.Content.Find.Execute("93-49=44", $true, $false, $false, $false, $false, $true, 1, $false, "53+15=68", 2) | Out-Null
$d.Content.Find.Execute("84-57=27", $true, $false, $false, $false, $false, $true, 1, $false, "51+46=97", 2) | Out-Null
$d.Content.Find.Execute("69+4=73", $true, $false, $false, $false, $false, $true, 1, $false, "95-1=94", 2) | Out-Null
$d.Content.Find.Execute("30-11=19", $true, $false, $false, $false, $false, $true, 1, $false, "37-21=16", 2) | Out-Null
$d.Content.Find.Execute("92-27=65", $true, $false, $false, $false, $false, $true, 1, $false, "48+2=50", 2) | Out-Null
$d.Content.Find.Execute("43+48=91", $true, $false, $false, $false, $false, $true, 1, $false, "45+14=59", 2) | Out-Null
$d.Content.Find.Execute("43-16=27", $true, $false, $false, $false, $false, $true, 1, $false, "42-41=1", 2) | Out-Null
$d.Content.Find.Execute("21+31=52", $true, $false, $false, $false, $false, $true, 1, $false, "31+68=99", 2) | Out-Null
$d.Content.Find.Execute("5+48=53", $true, $false, $false, $false, $false, $true, 1, $false, "64-23=41", 2) | Out-Null
$d.Content.Find.Execute("67+17=84", $true, $false, $false, $false, $false, $true, 1, $false, "91+2=93", 2) | Out-Null
$d.Content.Find.Execute("97-49=48", $true, $false, $false, $false, $false, $true, 1, $false, "91-81=10", 2) | Out-Null
$d.Content.Find.Execute("9+54=63", $true, $false, $false, $false, $false, $true, 1, $false, "73-46=27", 2) | Out-Null
$d.Content.Find.Execute("10+10=20", $true, $false, $false, $false, $false, $true, 1, $false, "34+32=66", 2) | Out-Null
$d.Content.Find.Execute("38-7=31", $true, $false, $false, $false, $false, $true, 1, $false, "30+52=82", 2) | Out-Null
$d.Content.Find.Execute("1+38=39", $true, $false, $false, $false, $false, $true, 1, $false, "97-74=23", 2) | Out-Null
$d.Content.Find.Execute("98-42=56", $true, $false, $false, $false, $false, $true, 1, $false, "69-50=19", 2) | Out-Null
$d.Content.Find.Execute("82-77=5", $true, $false, $false, $false, $false, $true, 1, $false, "24+13=37", 2) | Out-Null
$d.Content.Find.Execute("40+41=81", $true, $false, $false, $false, $false, $true, 1, $false, "74-20=54", 2) | Out-Null
$d.Content.Find.Execute("0+19=19", $true, $false, $false, $false, $false, $true, 1, $false, "20+22=42", 2) | Out-Null
$d.Content.Find.Execute("43+23=66", $true, $false, $false, $false, $false, $true, 1, $false, "63+33=96", 2) | Out-Null
$d.Content.Find.Execute("37+38=75", $true, $false, $false, $false, $false, $true, 1, $false, "49+32=81", 2) | Out-Null
$d.Content.Find.Execute("23-4=19", $true, $false, $false, $false, $false, $true, 1, $false, "38+55=93", 2) | Out-Null
$d.Content.Find.Execute("61-9=52", $true, $false, $false, $false, $false, $true, 1, $false, "12+13=25", 2) | Out-Null
$d.Content.Find.Execute("38+41=79", $true, $false, $false, $false, $false, $true, 1, $false, "60+34=94", 2) | Out-Null
$d.Content.Find.Execute("91-63=28", $true, $false, $false, $false, $false, $true, 1, $false, "20+53=73", 2) | Out-Null
$d.Content.Find.Execute("8+27=35", $true, $false, $false, $false, $false, $true, 1, $false, "6+45=51", 2) | Out-Null
$d.Content.Find.Execute("18-14=4", $true, $false, $false, $false, $false, $true, 1, $false, "46+15=61", 2) | Out-Null
$d.Content.Find.Execute("40-0=40", $true, $false, $false, $false, $false, $true, 1, $false, "68-64=4", 2) | Out-Null
$d.Content.Find.Execute("44-33=11", $true, $false, $false, $false, $false, $true, 1, $false, "15+23=38", 2) | Out-Null
$d.Content.Find.Execute("15+53=68", $true, $false, $false, $false, $false, $true, 1, $false, "92-82=10", 2) | Out-Null
$d.Content.Find.Execute("92-44=48", $true, $false, $false, $false, $false, $true, 1, $false, "90+9=99", 2) | Out-Null
$d.Content.Find.Execute("43-8=35", $true, $false, $false, $false, $false, $true, 1, $false, "47+33=80", 2) | Out-Null
$d.Content.Find.Execute("12+44=56", $true, $false, $false, $false, $false, $true, 1, $false, "64-5=59", 2) | Out-Null
$d.Content.Find.Execute("19-16=3", $true, $false, $false, $false, $false, $true, 1, $false, "21+22=43", 2) | Out-Null
$d.Content.Find.Execute("44+53=97", $true, $false, $false, $false, $false, $true, 1, $false, "6+74=80", 2) | Out-Null
$d.Content.Find.Execute("88-76=12", $true, $false, $false, $false, $false, $true, 1, $false, "73-31=42", 2) | Out-Null
$d.Content.Find.Execute("41+13=54", $true, $false, $false, $false, $false, $true, 1, $false, "72+16=88", 2) | Out-Null
$d.Content.Find.Execute("52-24=28", $true, $false, $false, $false, $false, $true, 1, $false, "19-13=6", 2) | Out-Null
$d.Content.Find.Execute("1+94=95", $true, $false, $false, $false, $false, $true, 1, $false, "57-40=17", 2) | Out-Null
$d.Content.Find.Execute("52-31=21", $true, $false, $false, $false, $false, $true, 1, $false, "45-38=7", 2) | Out-Null
$d.Content.Find.Execute("8+60=68", $true, $false, $false, $false, $false, $true, 1, $false, "42+33=75", 2) | Out-Null
$d.Content.Find.Execute("36+40=76", $true, $false, $false, $false, $false, $true, 1, $false, "66+28=94", 2) | Out-Null
$d.Content.Find.Execute("92-16=76", $true, $false, $false, $false, $false, $true, 1, $false, "87-74=13", 2) | Out-Null
$d.Content.Find.Execute("91+7=98", $true, $false, $false, $false, $false, $true, 1, $false, "18+67=85", 2) | Out-Null
$d.Content.Find.Execute("68-21=47", $true, $false, $false, $false, $false, $true, 1, $false, "72-5=67", 2) | Out-Null
$d.Content.Find.Execute("70-12=58", $true, $false, $false, $false, $false, $true, 1, $false, "58-14=44", 2) | Out-Null
$d.Content.Find.Execute("73+21=94", $true, $false, $false, $false, $false, $true, 1, $false, "33-17=16", 2) | Out-Null
$d.Content.Find.Execute("93-92=1", $true, $false, $false, $false, $false, $true, 1, $false, "66-15=51", 2) | Out-Null
$d.Content.Find.Execute("45+49=94", $true, $false, $false, $false, $false, $true, 1, $false, "45-27=18", 2) | Out-Null
$d.Content.Find.Execute("39+31=70", $true, $false, $false, $false, $false, $true, 1, $false, "26+11=37", 2) | Out-Null
$d.Content.Find.Execute("88-72=16", $true, $false, $false, $false, $false, $true, 1, $false, "17+34=51", 2) | Out-Null
$d.Content.Find.Execute("33+45=78", $true, $false, $false, $false, $false, $true, 1, $false, "56+3=59", 2) | Out-Null
$d.Content.Find.Execute("9+23=32", $true, $false, $false, $false, $false, $true, 1, $false, "74-69=5", 2) | Out-Null
$d.Content.Find.Execute("15+60=75", $true, $false, $false, $false, $false, $true, 1, $false, "54-36=18", 2) | Out-Null
$d.Content.Find.Execute("20+8=28", $true, $false, $false, $false, $false, $true, 1, $false, "35+26=61", 2) | Out-Null
$d.Content.Find.Execute("83+13=96", $true, $false, $false, $false, $false, $true, 1, $false, "28+30=58", 2) | Out-Null
$d.Content.Find.Execute("42+39=81", $true, $false, $false, $false, $false, $true, 1, $false, "59+28=87", 2) | Out-Null
$d.Content.Find.Execute("60-45=15", $true, $false, $false, $false, $false, $true, 1, $false, "50-48=2", 2) | Out-Null
$d.Content.Find.Execute("33+48=81", $true, $false, $false, $false, $false, $true, 1, $false, "29+2=31", 2) | Out-Null
$d.Content.Find.Execute("44-1=43", $true, $false, $false, $false, $false, $true, 1, $false, "28-17=11", 2) | Out-Null
$d.Content.Find.Execute("28+28=56", $true, $false, $false, $false, $false, $true, 1, $false, "40+59=99", 2) | Out-Null
$d.Content.Find.Execute("37+14=51", $true, $false, $false, $false, $false, $true, 1, $false, "31+46=77", 2) | Out-Null
$d.Content.Find.Execute("50+15=65", $true, $false, $false, $false, $false, $true, 1, $false, "44-11=33", 2) | Out-Null
$d.Content.Find.Execute("19+25=44", $true, $false, $false, $false, $false, $true, 1, $false, "11+10=21", 2) | Out-Null
$d.Content.Find.Execute("11+37=48", $true, $false, $false, $false, $false, $true, 1, $false, "53-26=27", 2) | Out-Null
$d.Content.Find.Execute("75-44=31", $true, $false, $false, $false, $false, $true, 1, $false, "76-7=69", 2) | Out-Null
$d.Content.Find.Execute("98-95=3", $true, $false, $false, $false, $false, $true, 1, $false, "57+30=87", 2) | Out-Null
$d.Content.Find.Execute("67-59=8", $true, $false, $false, $false, $false, $true, 1, $false, "16+82=98", 2) | Out-Null
$d.Content.Find.Execute("14+59=73", $true, $false, $false, $false, $false, $true, 1, $false, "78+10=88", 2) | Out-Null
$d.Content.Find.Execute("43-5=38", $true, $false, $false, $false, $false, $true, 1, $false, "50-21=29", 2) | Out-Null
$d.Content.Find.Execute("31+50=81", $true, $false, $false, $false, $false, $true, 1, $false, "75-48=27", 2) | Out-Null
$d.Content.Find.Execute("75-1=74", $true, $false, $false, $false, $false, $true, 1, $false, "85-48=37", 2) | Out-Null
$d.Content.Find.Execute("74+9=83", $true, $false, $false, $false, $false, $true, 1, $false, "32+31=63", 2) | Out-Null
$d.Content.Find.Execute("58-46=12", $true, $false, $false, $false, $false, $true, 1, $false, "11+15=26", 2) | Out-Null
$d.Content.Find.Execute("57+9=66", $true, $false, $false, $false, $false, $true, 1, $false, "45-39=6", 2) | Out-Null
$d.Content.Find.Execute("75-47=28", $true, $false, $false, $false, $false, $true, 1, $false, "60-20=40", 2) | Out-Null
$d.Content.Find.Execute("97-68=29", $true, $false, $false, $false, $false, $true, 1, $false, "5+24=29", 2) | Out-Null
$d.Content.Find.Execute("6+90=96", $true, $false, $false, $false, $false, $true, 1, $false, "17+36=53", 2) | Out-Null
$d.Content.Find.Execute("63-59=4", $true, $false, $false, $false, $false, $true, 1, $false, "5+63=68", 2) | Out-Null
$d.Content.Find.Execute("92-31=61", $true, $false, $false, $false, $false, $true, 1, $false, "75-10=65", 2) | Out-Null
$d.Content.Find.Execute("5+46=51", $true, $false, $false, $false, $false, $true, 1, $false, "99-10=89", 2) | Out-Null
$d.Content.Find.Execute("88-71=17", $true, $false, $false, $false, $false, $true, 1, $false, "57-14=43", 2) | Out-Null
$d.Content.Find.Execute("22+9=31", $true, $false, $false, $false, $false, $true, 1, $false, "99-26=73", 2) | Out-Null
$d.Content.Find.Execute("11+43=54", $true, $false, $false, $false, $false, $true, 1, $false, "34+8=42", 2) | Out-Null
$d.Content.Find.Execute("93-17=76", $true, $false, $false, $false, $false, $true, 1, $false, "73-24=49", 2) | Out-Null
$d.Content.Find.Execute("17-10=7", $true, $false, $false, $false, $false, $true, 1, $false, "77+16=93", 2) | Out-Null
$d.Content.Find.Execute("39+36=75", $true, $false, $false, $false, $false, $true, 1, $false, "88-47=41", 2) | Out-Null
$d.Content.Find.Execute("96-60=36", $true, $false, $false, $false, $false, $true, 1, $false, "63+14=77", 2) | Out-Null
$d.Content.Find.Execute("40+7=47", $true, $false, $false, $false, $false, $true, 1, $false, "42+23=65", 2) | Out-Null
$d.Content.Find.Execute("30-14=16", $true, $false, $false, $false, $false, $true, 1, $false, "2+69=71", 2) | Out-Null
$d.Content.Find.Execute("68-3=65", $true, $false, $false, $false, $false, $true, 1, $false, "30+23=53", 2) | Out-Null
$d.Content.Find.Execute("92-63=29", $true, $false, $false, $false, $false, $true, 1, $false, "75-74=1", 2) | Out-Null
$d.Content.Find.Execute("20-3=17", $true, $false, $false, $false, $false, $true, 1, $false, "3+81=84", 2) | Out-Null
$d.Content.Find.Execute("46+36=82", $true, $false, $false, $false, $false, $true, 1, $false, "41+31=72", 2) | Out-Null
$d.Content.Find.Execute("10+79=89", $true, $false, $false, $false, $false, $true, 1, $false, "73-33=40", 2) | Out-Null
$d.Content.Find.Execute("35+13=48", $true, $false, $false, $false, $false, $true, 1, $false, "52-25=27", 2) | Out-Null
$d.Content.Find.Execute("74-66=8", $true, $false, $false, $false, $false, $true, 1, $false, "19+32=51", 2) | Out-Null
$d.Content.Find.Execute("78+0=78", $true, $false, $false, $false, $false, $true, 1, $false, "38-31=7", 2) | Out-Null
$d.Content.Find.Execute("66-32=34", $true, $false, $false, $false, $false, $true, 1, $false, "15+84=99", 2) | Out-Null
